# feat: add 2022-Q4 data
#
# 1) Shift the "总计" (summary) sheet rows down by one and insert a new
#    "2022-Q4" row at the top of the data (row 2), renumbering the index
#    column (A) 0..7.
# 2) Insert a brand new worksheet named "2022-Q4" right after "总计" (i.e.
#    before the existing "2022-Q2" sheet) holding the per-fund holding
#    detail for the quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: shift "总计" rows down, insert the new 2022-Q4 summary row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Work from the bottom up so we don't clobber rows before they're copied.
$summary.Range("A8:D8").Copy($summary.Range("A9:D9"))
$summary.Range("A7:D7").Copy($summary.Range("A8:D8"))
$summary.Range("A6:D6").Copy($summary.Range("A7:D7"))
$summary.Range("A5:D5").Copy($summary.Range("A6:D6"))
$summary.Range("A4:D4").Copy($summary.Range("A5:D5"))
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# New row 2: 2022-Q4 data
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.21

# Renumber the index column (A) 0..7 top to bottom.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# Part 2: add the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$newSheet.Name = "2022-Q4"

# Match the sheetPr/pageMargins settings used by every other sheet in this
# workbook.
$newSheet.Outline.SummaryRow = [Microsoft.Office.Interop.Excel.XlSummaryRow]::xlSummaryBelow
$newSheet.Outline.SummaryColumn = [Microsoft.Office.Interop.Excel.XlSummaryColumn]::xlSummaryOnRight
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "003413", "华泰柏瑞新经济沪港深混合", "1.44", "94.26", "8.21", "0.1182", 1),
    @(1, "011355", "华泰柏瑞港股通时代机遇混合A", "0.70", "94.61", "6.30", "0.0441", 7),
    @(2, "460010", "华泰柏瑞亚洲领导企业混合（QDII）", "0.52", "97.17", "5.33", "0.0277", 9),
    @(3, "011356", "华泰柏瑞港股通时代机遇混合C", "0.39", "94.61", "6.30", "0.0246", 7)
)

# 基金代码 (fund code) and the text-like numeric columns (基金规模/股票总仓位/
# 仓位占比/持有市值) must stay text, matching the source data (leading zeros
# in fund codes, fixed decimal formatting), so force the format before
# assigning values.
$newSheet.Range("B2:G5").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Header + index-column styling to match the other sheets (bold, centered,
# bordered) - copy direct formatting only, values are left untouched.
$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$summary.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($fmt)
$summary.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial($fmt)
$excel.CutCopyMode = $false
